$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 8.6
$ws.Range("E2").Value = 57.8
$ws.Range("F2").Value = 10.4
$ws.Range("H2").Value = 56
$ws.Range("K2").Value = 57.7
$ws.Range("N2").Value = 52.28493729186943

# Row 3 updates
$ws.Range("D3").Value = 15.26
$ws.Range("E3").Value = 56.4
$ws.Range("F3").Value = 5.75
$ws.Range("J3").Value = 66
$ws.Range("K3").Value = 57.7
$ws.Range("N3").Value = 52.28493729186943
